$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each block below is a cyclic rotation of betting-odds rows (columns B:AD).
# Row order within the sheet (column A) stays put; the match/odds data that
# used to live in one row now belongs to the next row in the cycle.
$blocks = @(
  @(38, 39),
  @(45, 46),
  @(98, 99),
  @(105, 106, 107),
  @(114, 115),
  @(132, 133),
  @(145, 147),
  @(158, 160),
  @(164, 165),
  @(185, 186),
  @(211, 212),
  @(225, 226, 227),
  @(281, 282),
  @(376, 377, 378, 379, 380, 381)
)

foreach ($block in $blocks) {
  $data = @{}
  foreach ($r in $block) {
    $vals = @()
    for ($c = 2; $c -le 30; $c++) { $vals += $ws.Cells.Item($r, $c).Value2 }
    $data[$r] = $vals
  }
  $n = $block.Length
  for ($i = 0; $i -lt $n; $i++) {
    $dest = $block[$i]
    $src = $block[($i - 1 + $n) % $n]
    $vals = $data[$src]
    for ($c = 2; $c -le 30; $c++) { $ws.Cells.Item($dest, $c).Value = $vals[$c - 2] }
  }
}